# Updates cryptos list: refresh Price/Volume(1h) figures and shift in a new
# 'LEO' row at position 24 (pushing later rows down one slot; VeChain falls
# off the bottom of the fixed 50-row table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.026.36'
$ws.Range("E2").Value = '  -3.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.247.51'
$ws.Range("E3").Value = '  -4.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.88%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '496.11'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.37'
$ws.Range("E6").Value = '  -1.01%  '
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.292.50'
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0942'
$ws.Range("E10").Value = '  -3.77%  '
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("E12").Value = '  +0.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.63'
$ws.Range("E13").Value = '  -3.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.702.32'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.65'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '54.311.63'
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.293.84'
$ws.Range("E18").Value = '  -2.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.91'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("E20").Value = '  +1.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '303.97'
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.46'
$ws.Range("E22").Value = '  +4.23%  '
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.38'
$ws.Range("E24").Value = '  -2.79%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.80'
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.375'
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.399.07'
$ws.Range("E28").Value = '  -2.09%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.149'
$ws.Range("E29").Value = '  +2.04%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.14'
$ws.Range("E30").Value = '  +1.11%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.43'
$ws.Range("E31").Value = '  -3.67%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.60'
$ws.Range("E32").Value = '  -1.88%  '
$ws.Range("B33").Value = 'PEPE'
$ws.Range("C33").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0683'
$ws.Range("E33").Value = '  -2.55%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("E34").Value = '  +2.09%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.990'
$ws.Range("E36").Value = '  -0.68%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.07'
$ws.Range("E37").Value = '  +1.23%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.63'
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.19'
$ws.Range("E39").Value = '  +2.16%  '
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.874'
$ws.Range("E40").Value = '  +6.36%  '
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.64'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.44'
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.376'
$ws.Range("E43").Value = '  +1.59%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.40'
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.35'
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '126.22'
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.81'
$ws.Range("E47").Value = '  +2.88%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0891'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.546'
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '237.74'
$ws.Range("E50").Value = '  +0.47%  '
$ws.Range("B51").Value = 'Hedera'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0481'
$ws.Range("E51").Value = '  +1.60%  '
